$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-133: remaining 08.11.2025 entries (1-96) followed by
# new 09.11.2025 entries (1-29). The previous 07.11.2025 day (old rows 2-104)
# is dropped, and the final old row 134 (extra 08.11.2025 #30 duplicate) is
# removed since the sheet now ends at row 133.
$rows = @(
  @(45969,8.191000000000001,0.047,1,"08.11.20251"),
  @(45969.01041666666,6.212,0.498,2,"08.11.20252"),
  @(45969.02083333334,9.334,0,3,"08.11.20253"),
  @(45969.03125,8.616,0,4,"08.11.20254"),
  @(45969.04166666666,11.112,0,5,"08.11.20255"),
  @(45969.05208333334,13.144,0,6,"08.11.20256"),
  @(45969.0625,3.548,1.463,7,"08.11.20257"),
  @(45969.07291666666,1.428,3.758,8,"08.11.20258"),
  @(45969.08333333334,0.094,2.057,9,"08.11.20259"),
  @(45969.09375,0.019,5.8,10,"08.11.202510"),
  @(45969.10416666666,1.107,1.83,11,"08.11.202511"),
  @(45969.11458333334,0.731,5.566,12,"08.11.202512"),
  @(45969.125,0.299,4.991,13,"08.11.202513"),
  @(45969.13541666666,0,9.164,14,"08.11.202514"),
  @(45969.14583333334,0,6.924,15,"08.11.202515"),
  @(45969.15625,0,8.968999999999999,16,"08.11.202516"),
  @(45969.16666666666,0,19.656,17,"08.11.202517"),
  @(45969.17708333334,0,14.908,18,"08.11.202518"),
  @(45969.1875,0,5.321,19,"08.11.202519"),
  @(45969.19791666666,0,5.875,20,"08.11.202520"),
  @(45969.20833333334,0,33.339,21,"08.11.202521"),
  @(45969.21875,0,42.556,22,"08.11.202522"),
  @(45969.22916666666,0,38.377,23,"08.11.202523"),
  @(45969.23958333334,0,13.999,24,"08.11.202524"),
  @(45969.25,0,60.406,25,"08.11.202525"),
  @(45969.26041666666,0,20.337,26,"08.11.202526"),
  @(45969.27083333334,0,13.395,27,"08.11.202527"),
  @(45969.28125,0,12.136,28,"08.11.202528"),
  @(45969.29166666666,0,7.366,29,"08.11.202529"),
  @(45969.30208333334,1.866,16.943,30,"08.11.202530"),
  @(45969.3125,1.204,3.135,31,"08.11.202531"),
  @(45969.32291666666,0.864,12.297,32,"08.11.202532"),
  @(45969.33333333334,21.532,0.481,33,"08.11.202533"),
  @(45969.34375,46.86,0,34,"08.11.202534"),
  @(45969.35416666666,21.981,0,35,"08.11.202535"),
  @(45969.36458333334,28.526,0,36,"08.11.202536"),
  @(45969.375,3.272,0,37,"08.11.202537"),
  @(45969.38541666666,3.48,0,38,"08.11.202538"),
  @(45969.39583333334,14.001,0,39,"08.11.202539"),
  @(45969.40625,13.29,0,40,"08.11.202540"),
  @(45969.41666666666,1.676,22.988,41,"08.11.202541"),
  @(45969.42708333334,0,19.586,42,"08.11.202542"),
  @(45969.4375,0.052,13.749,43,"08.11.202543"),
  @(45969.44791666666,3.04,3.246,44,"08.11.202544"),
  @(45969.45833333334,0,19.69,45,"08.11.202545"),
  @(45969.46875,0,19.656,46,"08.11.202546"),
  @(45969.47916666666,0.954,1.067,47,"08.11.202547"),
  @(45969.48958333334,1.945,3.585,48,"08.11.202548"),
  @(45969.5,1.077,2.395,49,"08.11.202549"),
  @(45969.51041666666,1.007,3.997,50,"08.11.202550"),
  @(45969.52083333334,7.823,0.001,51,"08.11.202551"),
  @(45969.53125,8.351000000000001,0.038,52,"08.11.202552"),
  @(45969.54166666666,4.253,0.896,53,"08.11.202553"),
  @(45969.55208333334,3.044,2.062,54,"08.11.202554"),
  @(45969.5625,5.368,0.471,55,"08.11.202555"),
  @(45969.57291666666,4.974,0.762,56,"08.11.202556"),
  @(45969.58333333334,3.129,1.833,57,"08.11.202557"),
  @(45969.59375,0.144,15.945,58,"08.11.202558"),
  @(45969.60416666666,0,18.451,59,"08.11.202559"),
  @(45969.61458333334,0,25.108,60,"08.11.202560"),
  @(45969.625,0,37.958,61,"08.11.202561"),
  @(45969.63541666666,0,46.902,62,"08.11.202562"),
  @(45969.64583333334,0,50.968,63,"08.11.202563"),
  @(45969.65625,0.026,12.468,64,"08.11.202564"),
  @(45969.66666666666,0.128,8.068,65,"08.11.202565"),
  @(45969.67708333334,0,70.008,66,"08.11.202566"),
  @(45969.6875,0,101.736,67,"08.11.202567"),
  @(45969.69791666666,0.064,35.8,68,"08.11.202568"),
  @(45969.70833333334,0,32.071,69,"08.11.202569"),
  @(45969.71875,0,44.327,70,"08.11.202570"),
  @(45969.72916666666,0,59.777,71,"08.11.202571"),
  @(45969.73958333334,0,40.255,72,"08.11.202572"),
  @(45969.75,0.827,12.549,73,"08.11.202573"),
  @(45969.76041666666,0,49.195,74,"08.11.202574"),
  @(45969.77083333334,0,39.369,75,"08.11.202575"),
  @(45969.78125,0.03,11.437,76,"08.11.202576"),
  @(45969.79166666666,3.42,0.204,77,"08.11.202577"),
  @(45969.80208333334,6.567,0.089,78,"08.11.202578"),
  @(45969.8125,0,4.267,79,"08.11.202579"),
  @(45969.82291666666,0,2.275,80,"08.11.202580"),
  @(45969.83333333334,0,18.865,81,"08.11.202581"),
  @(45969.84375,0,6.799,82,"08.11.202582"),
  @(45969.85416666666,0.598,1.222,83,"08.11.202583"),
  @(45969.86458333334,6.281,0.002,84,"08.11.202584"),
  @(45969.875,37.09,0,85,"08.11.202585"),
  @(45969.88541666666,57.436,0,86,"08.11.202586"),
  @(45969.89583333334,59.553,0,87,"08.11.202587"),
  @(45969.90625,25.54,0,88,"08.11.202588"),
  @(45969.91666666666,63.154,0,89,"08.11.202589"),
  @(45969.92708333334,62.835,0,90,"08.11.202590"),
  @(45969.92708333334,62.835,0,90,"08.11.202590"),
  @(45969.9375,34.512,0,91,"08.11.202591"),
  @(45969.9375,34.512,0,91,"08.11.202591"),
  @(45969.94791666666,0.051,41.195,92,"08.11.202592"),
  @(45969.94791666666,0.051,41.195,92,"08.11.202592"),
  @(45969.95833333334,0,36.848,93,"08.11.202593"),
  @(45969.95833333334,0,36.848,93,"08.11.202593"),
  @(45969.96875,0,28.94,94,"08.11.202594"),
  @(45969.96875,0,28.94,94,"08.11.202594"),
  @(45969.97916666666,0,23.497,95,"08.11.202595"),
  @(45969.97916666666,0,23.497,95,"08.11.202595"),
  @(45969.98958333334,0,22.461,96,"08.11.202596"),
  @(45969.98958333334,0,22.461,96,"08.11.202596"),
  @(45970,0,7.583,1,"09.11.20251"),
  @(45970.01041666666,0.12,11.998,2,"09.11.20252"),
  @(45970.02083333334,1.23,1.44,3,"09.11.20253"),
  @(45970.03125,3.624,3.678,4,"09.11.20254"),
  @(45970.04166666666,3.204,0.61,5,"09.11.20255"),
  @(45970.05208333334,2.952,0.5600000000000001,6,"09.11.20256"),
  @(45970.0625,2.939,1.127,7,"09.11.20257"),
  @(45970.07291666666,0.145,3.377,8,"09.11.20258"),
  @(45970.08333333334,0,43.665,9,"09.11.20259"),
  @(45970.09375,0.916,9.532,10,"09.11.202510"),
  @(45970.10416666666,0.028,5.968,11,"09.11.202511"),
  @(45970.11458333334,0,22.399,12,"09.11.202512"),
  @(45970.125,0,10.42,13,"09.11.202513"),
  @(45970.13541666666,0,13.868,14,"09.11.202514"),
  @(45970.14583333334,0,22.114,15,"09.11.202515"),
  @(45970.15625,0,39.259,16,"09.11.202516"),
  @(45970.16666666666,0,42.196,17,"09.11.202517"),
  @(45970.17708333334,0,57.057,18,"09.11.202518"),
  @(45970.1875,0,42.629,19,"09.11.202519"),
  @(45970.19791666666,0,26.289,20,"09.11.202520"),
  @(45970.20833333334,0,34.201,21,"09.11.202521"),
  @(45970.21875,0,36.947,22,"09.11.202522"),
  @(45970.22916666666,0,30.282,23,"09.11.202523"),
  @(45970.23958333334,0,36.687,24,"09.11.202524"),
  @(45970.25,0,49.486,25,"09.11.202525"),
  @(45970.26041666666,0,45.473,26,"09.11.202526"),
  @(45970.27083333334,0,24.09,27,"09.11.202527"),
  @(45970.28125,0.499,4.425,28,"09.11.202528"),
  @(45970.29166666666,0,11.828,29,"09.11.202529")
)

$numRows = $rows.Count
$numCols = 5
$data = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}

$ws.Range("A2:E133").Value = $data

# The old sheet had 133 data rows (ending at row 134); the refreshed data
# only needs 132 data rows (ending at row 133), so drop the now-unused row.
$ws.Rows(134).Delete()
